$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.60230016708374
$ws.Range("B1").Value = 4.174076080322266
$ws.Range("C1").Value = 3.308237314224243
$ws.Range("D1").Value = 2.019135475158691
$ws.Range("E1").Value = 0.7408944368362427
